$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# D-column values are forced to text format before assignment so that
# numeric-looking strings (e.g. "218.18") are not auto-converted to numbers
# by Excel, preserving the exact text representation used by the site.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.132.09"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.21"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.18"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5285"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2607"
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.44"
$ws.Range("E10").Value = "  -2.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07781"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.513"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.656.45"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8222"
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.49"
$ws.Range("E16").Value = "  +1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.151.53"
$ws.Range("E17").Value = "  -0.51%  "
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.579"
$ws.Range("E19").Value = "  -1.93%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.06"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.10"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.054"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "141.34"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.289"
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.21"
$ws.Range("E27").Value = "  +0.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.442"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05942"
$ws.Range("E29").Value = "  -2.99%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.526"
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.582"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9543"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.789"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("E36").Value = "  -0.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5713"
$ws.Range("E37").Value = "  +0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01620"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.812"
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8476"
$ws.Range("E40").Value = "  -0.95%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.06"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.024.41"
$ws.Range("E43").Value = "  +0.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.802.56"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.008"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.492"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4297"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05154"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.815"
$ws.Range("E50").Value = "  -3.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09720"
$ws.Range("E51").Value = "  -0.17%  "
